$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 778
$ws.Range("I12").Value = 778
$ws.Range("K12").Value = 778
$ws.Range("M12").Value = -608
$ws.Range("H18").Value = 10950
$ws.Range("I18").Value = 10950
$ws.Range("K18").Value = 10950
$ws.Range("M18").Value = -10666
$ws.Range("H28").Value = 1220.1666
$ws.Range("I28").Value = 1541
$ws.Range("J28").Value = 257.66666
$ws.Range("K28").Value = 1541
$ws.Range("L28").Value = 257.66666
$ws.Range("M28").Value = -1056
$ws.Range("N28").Value = -1227.66666
$ws.Range("H62").Value = 1780.7142
$ws.Range("J62").Value = 1761.3334
$ws.Range("L62").Value = 1761.3334
$ws.Range("N62").Value = -3009.3334
$ws.Range("H65").Value = 1780.7142
$ws.Range("J65").Value = 1761.3334
$ws.Range("L65").Value = 8806.666999999999
$ws.Range("N65").Value = -15046.667
$ws.Range("H125").Value = 11083.286
$ws.Range("I125").Value = 1750
$ws.Range("J125").Value = 14816.6
$ws.Range("K125").Value = 15750
$ws.Range("L125").Value = 133349.4
$ws.Range("M125").Value = -13290
$ws.Range("N125").Value = -138269.4
$ws.Range("H132").Value = 2096.7778
$ws.Range("I132").Value = 2096.7778
$ws.Range("K132").Value = 6290.3334
$ws.Range("M132").Value = -3760.3334
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1296.125
$ws.Range("I45").Value = 1296.125
$ws.Range("K45").Value = 1296.125
$ws.Range("M45").Value = -919.125
$ws.Range("H74").Value = 4558.1665
$ws.Range("I74").Value = 3669.8
$ws.Range("K74").Value = 3669.8
$ws.Range("M74").Value = -2795.8
$ws.Range("H77").Value = 4558.1665
$ws.Range("I77").Value = 3669.8
$ws.Range("K77").Value = 18349
$ws.Range("M77").Value = -13981
$ws.Range("H88").Value = 2573.6667
$ws.Range("I88").Value = 2652
$ws.Range("J88").Value = 2495.3333
$ws.Range("K88").Value = 2652
$ws.Range("L88").Value = 2495.3333
$ws.Range("M88").Value = -2246
$ws.Range("N88").Value = -3307.3333
$ws.Range("H91").Value = 2573.6667
$ws.Range("I91").Value = 2652
$ws.Range("J91").Value = 2495.3333
$ws.Range("K91").Value = 2652
$ws.Range("L91").Value = 2495.3333
$ws.Range("M91").Value = -1248
$ws.Range("N91").Value = -5303.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1449.2
$ws.Range("I86").Value = 1395
$ws.Range("J86").Value = 1666
$ws.Range("K86").Value = 1395
$ws.Range("L86").Value = 1666
$ws.Range("M86").Value = -272
$ws.Range("N86").Value = -3912
$ws.Range("H89").Value = 1449.2
$ws.Range("I89").Value = 1395
$ws.Range("J89").Value = 1666
$ws.Range("K89").Value = 6975
$ws.Range("L89").Value = 8330
$ws.Range("M89").Value = -1359
$ws.Range("N89").Value = -19562

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1811.8334
$ws.Range("I31").Value = 774.25
$ws.Range("K31").Value = 774.25
$ws.Range("M31").Value = -479.25
$ws.Range("H34").Value = 1811.8334
$ws.Range("I34").Value = 774.25
$ws.Range("K34").Value = 774.25
$ws.Range("M34").Value = -572.25
$ws.Range("H105").Value = 2543.5557
$ws.Range("I105").Value = 1495.4445
$ws.Range("K105").Value = 1495.4445
$ws.Range("M105").Value = 251.5554999999999
$ws.Range("H107").Value = 722.5
$ws.Range("I107").Value = 290
$ws.Range("K107").Value = 290
$ws.Range("M107").Value = 1630

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 477
$ws.Range("I44").Value = 638
$ws.Range("K44").Value = 1914
$ws.Range("M44").Value = -1516
$ws.Range("H69").Value = 973.75
$ws.Range("J69").Value = 965.3333
$ws.Range("L69").Value = 2895.9999
$ws.Range("N69").Value = -4517.9999
$ws.Range("H72").Value = 973.75
$ws.Range("J72").Value = 965.3333
$ws.Range("L72").Value = 8687.9997
$ws.Range("N72").Value = -16799.9997
$ws.Range("H81").Value = 2153.25
$ws.Range("J81").Value = 2450
$ws.Range("L81").Value = 7350
$ws.Range("N81").Value = -9596
$ws.Range("H84").Value = 2153.25
$ws.Range("J84").Value = 2450
$ws.Range("L84").Value = 22050
$ws.Range("N84").Value = -33282
$ws.Range("H92").Value = 423.8
$ws.Range("J92").Value = 439.66666
$ws.Range("L92").Value = 1318.99998
$ws.Range("N92").Value = -3814.99998
$ws.Range("H109").Value = 962.9
$ws.Range("I109").Value = 692.1111
$ws.Range("J109").Value = 3400
$ws.Range("K109").Value = 2076.3333
$ws.Range("L109").Value = 10200
$ws.Range("M109").Value = -1036.3333
$ws.Range("N109").Value = -12280
$ws.Range("H137").Value = 3998.5
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 3998.5
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 11995.5
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -22195.5
$ws.Range("H140").Value = 8450.362999999999
$ws.Range("I140").Value = 3986.6667
$ws.Range("K140").Value = 11960.0001
$ws.Range("M140").Value = -6780.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 48999.5
$ws.Range("J86").Value = 48999.5
$ws.Range("L86").Value = 48999.5
$ws.Range("N86").Value = -51371.5
$ws.Range("H89").Value = 48999.5
$ws.Range("J89").Value = 48999.5
$ws.Range("L89").Value = 146998.5
$ws.Range("N89").Value = -158854.5
$ws.Range("H102").Value = 2373
$ws.Range("I102").Value = 2431.2222
$ws.Range("K102").Value = 2431.2222
$ws.Range("M102").Value = -809.2222000000002
$ws.Range("H107").Value = 3311.875
$ws.Range("J107").Value = 3571.4285
$ws.Range("L107").Value = 3571.4285
$ws.Range("N107").Value = -7411.4285
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 706.8
$ws.Range("I122").Value = 946.5
$ws.Range("J122").Value = 347.25
$ws.Range("K122").Value = 2839.5
$ws.Range("L122").Value = 1041.75
$ws.Range("M122").Value = -389.5
$ws.Range("N122").Value = -5941.75
$ws.Range("H126").Value = 2165.3333
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 2330.6667
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 6992.000100000001
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -11932.0001
$ws.Range("H132").Value = 19990
$ws.Range("I132").Value = 19990
$ws.Range("K132").Value = 59970
$ws.Range("M132").Value = -57440

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 6794
$ws.Range("I16").Value = 725.3333
$ws.Range("K16").Value = 725.3333
$ws.Range("M16").Value = -555.3333
$ws.Range("H22").Value = 654
$ws.Range("I22").Value = 627.6667
$ws.Range("K22").Value = 627.6667
$ws.Range("M22").Value = -332.6667
$ws.Range("H27").Value = 654
$ws.Range("I27").Value = 627.6667
$ws.Range("K27").Value = 627.6667
$ws.Range("M27").Value = -520.6667
$ws.Range("H40").Value = 3798.3333
$ws.Range("I40").Value = 2697.5
$ws.Range("K40").Value = 2697.5
$ws.Range("M40").Value = -2561.5
$ws.Range("H46").Value = 1696.6666
$ws.Range("I46").Value = 350.22223
$ws.Range("K46").Value = 350.22223
$ws.Range("M46").Value = -162.22223
$ws.Range("H68").Value = 2149.6667
$ws.Range("J68").Value = 2224.5
$ws.Range("L68").Value = 2224.5
$ws.Range("N68").Value = -3722.5
$ws.Range("H71").Value = 2149.6667
$ws.Range("J71").Value = 2224.5
$ws.Range("L71").Value = 11122.5
$ws.Range("N71").Value = -18610.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 26465.666
$ws.Range("J56").Value = 26465.666
$ws.Range("L56").Value = 26465.666
$ws.Range("N56").Value = -27893.666
